# Auto-generated Excel COM-interop edit script
# Commit: Update IESO report from GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report generation timestamp
$ws.Range("A1").Value = "CreatedAt: 2026-01-20T17:07:33"

# Refresh predispatch hourly LMP figures for hours 19-24 (columns U-Z)
$ws.Range("U4").Value = 85.95999999999999
$ws.Range("V4").Value = 99
$ws.Range("W4").Value = 58.32
$ws.Range("X4").Value = 63.03
$ws.Range("Y4").Value = 49.56
$ws.Range("Z4").Value = 51.07
$ws.Range("U6").Value = -2.41
$ws.Range("V6").Value = -3.27
$ws.Range("W6").Value = -1.63
$ws.Range("X6").Value = -1.13
$ws.Range("Y6").Value = -0.55
$ws.Range("Z6").Value = 0
$ws.Range("U9").Value = 84.31999999999999
$ws.Range("V9").Value = 97.03
$ws.Range("W9").Value = 58.03
$ws.Range("X9").Value = 63.66
$ws.Range("Y9").Value = 50.92
$ws.Range("Z9").Value = 52.27
$ws.Range("U11").Value = -4.05
$ws.Range("V11").Value = -5.24
$ws.Range("W11").Value = -1.92
$ws.Range("X11").Value = -0.51
$ws.Range("Y11").Value = 0.8100000000000001
$ws.Range("Z11").Value = 1.2
$ws.Range("U14").Value = 84.40000000000001
$ws.Range("V14").Value = 104.18
$ws.Range("W14").Value = 58.03
$ws.Range("X14").Value = 63.66
$ws.Range("Y14").Value = 50.92
$ws.Range("Z14").Value = 61.93
$ws.Range("V15").Value = 7.15
$ws.Range("Z15").Value = 9.66
$ws.Range("U16").Value = -3.97
$ws.Range("V16").Value = -5.24
$ws.Range("W16").Value = -1.92
$ws.Range("X16").Value = -0.51
$ws.Range("Y16").Value = 0.8100000000000001
$ws.Range("Z16").Value = 1.2
$ws.Range("U19").Value = 84.56999999999999
$ws.Range("V19").Value = 97.48999999999999
$ws.Range("W19").Value = 57.42
$ws.Range("X19").Value = 62.12
$ws.Range("Y19").Value = 48.97
$ws.Range("Z19").Value = 50.66
$ws.Range("U21").Value = -3.81
$ws.Range("V21").Value = -4.78
$ws.Range("W21").Value = -2.53
$ws.Range("X21").Value = -2.05
$ws.Range("Y21").Value = -1.13
$ws.Range("Z21").Value = -0.41
$ws.Range("U24").Value = 84.56999999999999
$ws.Range("V24").Value = 97.48999999999999
$ws.Range("W24").Value = 57.42
$ws.Range("X24").Value = 62.12
$ws.Range("Y24").Value = 48.97
$ws.Range("Z24").Value = 50.66
$ws.Range("U26").Value = -3.81
$ws.Range("V26").Value = -4.78
$ws.Range("W26").Value = -2.53
$ws.Range("X26").Value = -2.05
$ws.Range("Y26").Value = -1.13
$ws.Range("Z26").Value = -0.41
$ws.Range("U29").Value = 82.90000000000001
$ws.Range("V29").Value = 95.67
$ws.Range("W29").Value = 56.29
$ws.Range("X29").Value = 60.88
$ws.Range("Y29").Value = 48.13
$ws.Range("Z29").Value = 50.02
$ws.Range("U31").Value = -5.47
$ws.Range("V31").Value = -6.6
$ws.Range("W31").Value = -3.66
$ws.Range("X31").Value = -3.29
$ws.Range("Y31").Value = -1.97
$ws.Range("Z31").Value = -1.05
$ws.Range("U34").Value = 85.3
$ws.Range("V34").Value = 103.22
$ws.Range("W34").Value = 59.06
$ws.Range("X34").Value = 65.41
$ws.Range("Y34").Value = 52.68
$ws.Range("Z34").Value = 61.31
$ws.Range("V35").Value = 5.64
$ws.Range("Z35").Value = 7.61
$ws.Range("U36").Value = -3.07
$ws.Range("V36").Value = -4.68
$ws.Range("W36").Value = -0.89
$ws.Range("X36").Value = 1.24
$ws.Range("Y36").Value = 2.58
$ws.Range("U39").Value = 85.95999999999999
$ws.Range("V39").Value = 99
$ws.Range("W39").Value = 58.32
$ws.Range("X39").Value = 63.03
$ws.Range("Y39").Value = 49.56
$ws.Range("Z39").Value = 51.07
$ws.Range("U41").Value = -2.41
$ws.Range("V41").Value = -3.27
$ws.Range("W41").Value = -1.63
$ws.Range("X41").Value = -1.13
$ws.Range("Y41").Value = -0.55
$ws.Range("Z41").Value = 0
$ws.Range("U44").Value = 95.54000000000001
$ws.Range("V44").Value = 110.32
$ws.Range("W44").Value = 64.45999999999999
$ws.Range("X44").Value = 68.7
$ws.Range("Y44").Value = 53.76
$ws.Range("Z44").Value = 55.03
$ws.Range("U46").Value = 7.17
$ws.Range("V46").Value = 8.050000000000001
$ws.Range("W46").Value = 4.51
$ws.Range("X46").Value = 4.53
$ws.Range("Y46").Value = 3.66
$ws.Range("Z46").Value = 3.96
$ws.Range("U49").Value = 81.67
$ws.Range("V49").Value = 100.95
$ws.Range("W49").Value = 61.74
$ws.Range("X49").Value = 65.81
$ws.Range("Y49").Value = 51.07
$ws.Range("Z49").Value = 52.38
$ws.Range("U51").Value = -6.7
$ws.Range("V51").Value = -1.31
$ws.Range("W51").Value = 1.79
$ws.Range("X51").Value = 1.65
$ws.Range("Y51").Value = 0.97
$ws.Range("Z51").Value = 1.31
$ws.Range("U54").Value = 85.22
$ws.Range("V54").Value = 98.70999999999999
$ws.Range("W54").Value = 57.87
$ws.Range("X54").Value = 61.76
$ws.Range("Y54").Value = 48.55
$ws.Range("Z54").Value = 50.56
$ws.Range("U56").Value = -3.15
$ws.Range("V56").Value = -3.55
$ws.Range("W56").Value = -2.08
$ws.Range("X56").Value = -2.41
$ws.Range("Y56").Value = -1.55
$ws.Range("Z56").Value = -0.51
$ws.Range("U59").Value = 93.22
$ws.Range("V59").Value = 107.88
$ws.Range("W59").Value = 63.04
$ws.Range("X59").Value = 67.26000000000001
$ws.Range("Y59").Value = 52.63
$ws.Range("Z59").Value = 53.81
$ws.Range("U61").Value = 4.85
$ws.Range("V61").Value = 5.61
$ws.Range("W61").Value = 3.09
$ws.Range("X61").Value = 3.09
$ws.Range("Y61").Value = 2.53
$ws.Range("Z61").Value = 2.74
$ws.Range("U64").Value = 94.31
$ws.Range("V64").Value = 109.03
$ws.Range("W64").Value = 63.71
$ws.Range("X64").Value = 67.97
$ws.Range("Y64").Value = 53.19
$ws.Range("Z64").Value = 54.27
$ws.Range("U66").Value = 5.94
$ws.Range("V66").Value = 6.76
$ws.Range("W66").Value = 3.76
$ws.Range("X66").Value = 3.81
$ws.Range("Y66").Value = 3.08
$ws.Range("Z66").Value = 3.2
$ws.Range("V69").Value = 176
$ws.Range("W69").Value = 63.91
$ws.Range("X69").Value = 68.19
$ws.Range("Y69").Value = 53.41
$ws.Range("Z69").Value = 54.56
$ws.Range("U71").Value = 6.45
$ws.Range("V71").Value = 7.34
$ws.Range("W71").Value = 3.96
$ws.Range("X71").Value = 4.02
$ws.Range("Y71").Value = 3.31
$ws.Range("Z71").Value = 3.49
$ws.Range("U72").Value = 37.3
$ws.Range("V72").Value = 66.39
$ws.Range("U74").Value = 92.63
$ws.Range("V74").Value = 107.76
$ws.Range("W74").Value = 62.84
$ws.Range("X74").Value = 67.05
$ws.Range("Y74").Value = 52.57
$ws.Range("Z74").Value = 53.7
$ws.Range("U76").Value = 4.26
$ws.Range("V76").Value = 5.5
$ws.Range("W76").Value = 2.89
$ws.Range("X76").Value = 2.88
$ws.Range("Y76").Value = 2.47
$ws.Range("Z76").Value = 2.63
$ws.Range("U79").Value = 92.33
$ws.Range("V79").Value = 107.23
$ws.Range("W79").Value = 62.95
$ws.Range("X79").Value = 67.08
$ws.Range("Y79").Value = 52.74
$ws.Range("Z79").Value = 54.05
$ws.Range("U81").Value = 3.95
$ws.Range("V81").Value = 4.97
$ws.Range("W81").Value = 3.01
$ws.Range("X81").Value = 2.92
$ws.Range("Y81").Value = 2.64
$ws.Range("Z81").Value = 2.98
$ws.Range("U84").Value = 83.45
$ws.Range("V84").Value = 96.66
$ws.Range("W84").Value = 56.72
$ws.Range("X84").Value = 60.71
$ws.Range("Y84").Value = 47.72
$ws.Range("Z84").Value = 51.22
$ws.Range("U86").Value = -4.92
$ws.Range("V86").Value = -5.61
$ws.Range("W86").Value = -3.23
$ws.Range("X86").Value = -3.46
$ws.Range("Y86").Value = -2.39
$ws.Range("Z86").Value = 0.15
$ws.Range("U89").Value = 82.90000000000001
$ws.Range("V89").Value = 95.67
$ws.Range("W89").Value = 56.29
$ws.Range("X89").Value = 60.88
$ws.Range("Y89").Value = 48.08
$ws.Range("Z89").Value = 50.02
$ws.Range("U91").Value = -5.47
$ws.Range("V91").Value = -6.6
$ws.Range("W91").Value = -3.66
$ws.Range("X91").Value = -3.29
$ws.Range("Y91").Value = -2.02
$ws.Range("Z91").Value = -1.05
